# drill down system tab clockin
# Clear out the three stale test-data cells that are no longer needed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("A5").ClearContents()

$ws.Range("B8").Select()
